# UI Validations added V.66
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UIValidations")

# Correct TotalCharge value (row 33, column C) - remove the stale "1,800.00 USD::" prefix
$ws.Range("C33").Value = "1,057.50 USD"

# Correct Value2 value (row 43, column D) - remove the stale "NA::" prefix
$ws.Range("D43").Value = "351.80 USD"

# Update the view state: scroll so row 19 is the top row and select C52
$ws.Range("C52").Select()
$excel.ActiveWindow.ScrollRow = 19
